$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two control point values that remain
$ws.Range("B2").Value = 413
$ws.Range("B3").Value = 217

# Remove the last two rows (A4:B4 and A5:B5) entirely
$ws.Range("A4:B5").Delete()
